# "excel sheet1 test22 and sheet2added"
#
# 1. Sheet1!A2 changes from "test2" -> "test22"
# 2. A new worksheet "Sheet2" is inserted right after Sheet1, with
#    A1 = "sheet2_t1" and A2 = "sheet2_t2"
# 3. The selection on both sheets ends up on cell A3, and Sheet1 stays the
#    active/visible tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet immediately after Sheet1 (Add() with no args would
# insert it *before* the active sheet, so pass Sheet1 as the After target).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 first so the shared-string table picks up "sheet2_t1" /
# "sheet2_t2" before "test22".
$ws2.Range("A1").Value = "sheet2_t1"
$ws2.Range("A2").Value = "sheet2_t2"

# Update Sheet1's existing value.
$ws1.Range("A2").Value = "test22"

# Move the selection on each sheet to A3.
[void]$ws2.Range("A3").Select()
[void]$ws1.Activate()
[void]$ws1.Range("A3").Select()
